$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8102841940994949
$ws.Range("C2").Value = 0.2220612147349925
$ws.Range("E2").Value = 0.4644618158516849
$ws.Range("F2").Value = 2.059188574304642
$ws.Range("G2").Value = 0.5321724232187819
$ws.Range("H2").Value = 0.6450555677373586
$ws.Range("I2").Value = 0.4625974569064297
$ws.Range("J2").Value = 0.03585200113963216
$ws.Range("M2").Value = 0.5596912608547484
$ws.Range("B3").Value = 0.7121606348314344
$ws.Range("C3").Value = 0.1935999399594834
$ws.Range("E3").Value = 0.4612373741089968
$ws.Range("F3").Value = 2.047632612782522
$ws.Range("G3").Value = 0.5293149327334703
$ws.Range("H3").Value = 0.6498035695886983
$ws.Range("I3").Value = 0.4702532687446812
$ws.Range("J3").Value = 0.03667570645015683
$ws.Range("M3").Value = 0.5244027752425211
$ws.Range("B4").Value = 0.6518435782678864
$ws.Range("C4").Value = 0.176074140028561
$ws.Range("E4").Value = 0.4594086031030713
$ws.Range("F4").Value = 2.041930474746664
$ws.Range("G4").Value = 0.5282013462138337
$ws.Range("H4").Value = 0.6532225008564012
$ws.Range("I4").Value = 0.4754185898894256
$ws.Range("J4").Value = 0.03721014092745012
$ws.Range("M4").Value = 0.5028929364600998
$ws.Range("B5").Value = 0.6272473984886631
$ws.Range("C5").Value = 0.1689196198656759
$ws.Range("E5").Value = 0.4587014321102245
$ws.Range("F5").Value = 2.039956416641445
$ws.Range("G5").Value = 0.5279077145575144
$ws.Range("H5").Value = 0.6547420421762524
$ws.Range("I5").Value = 0.4776398726074582
$ws.Range("J5").Value = 0.03743512240291169
$ws.Range("M5").Value = 0.4941673922509864
$ws.Range("B6").Value = 0.6231622499318519
$ws.Range("C6").Value = 0.1677308577021677
$ws.Range("E6").Value = 0.4585863086612889
$ws.Range("F6").Value = 2.039649719802881
$ws.Range("G6").Value = 0.5278686014122798
$ws.Range("H6").Value = 0.6550019789152941
$ws.Range("I6").Value = 0.4780157300268186
$ws.Range("J6").Value = 0.03747291451266754
$ws.Range("M6").Value = 0.4927209405688444
$ws.Range("B7").Value = 0.6515119312614672
$ws.Range("C7").Value = 0.1759777026366578
$ws.Range("E7").Value = 0.4593989116975266
$ws.Range("F7").Value = 2.041902437339473
$ws.Range("G7").Value = 0.5281967390278766
$ws.Range("H7").Value = 0.6532424830041776
$ws.Range("I7").Value = 0.4754480763116753
$ws.Range("J7").Value = 0.03721314599437342
$ws.Range("M7").Value = 0.5027750988018695
$ws.Range("B8").Value = 0.776465919080124
$ws.Range("C8").Value = 0.2122582755017675
$ws.Range("E8").Value = 0.4633187272089216
$ws.Range("F8").Value = 2.054914450398456
$ws.Range("G8").Value = 0.5310535315287694
$ws.Range("H8").Value = 0.6465879400961256
$ws.Range("I8").Value = 0.4651404721211776
$ws.Range("J8").Value = 0.03613005359284971
$ws.Range("M8").Value = 0.5474912828639731
$ws.Range("B9").Value = 1.020931821948238
$ws.Range("C9").Value = 0.2830043815093575
$ws.Range("E9").Value = 0.4722010382391488
$ws.Range("F9").Value = 2.091523772365306
$ws.Range("G9").Value = 0.5417874620377887
$ws.Range("H9").Value = 0.6375505305093299
$ws.Range("I9").Value = 0.4486344665106792
$ws.Range("J9").Value = 0.03423432201524701
$ws.Range("M9").Value = 0.6364205811578643
$ws.Range("B10").Value = 1.200183030692415
$ws.Range("C10").Value = 0.3347437131761239
$ws.Range("E10").Value = 0.4794529662036808
$ws.Range("F10").Value = 2.125241982735446
$ws.Range("G10").Value = 0.5528679097414937
$ws.Range("H10").Value = 0.6333793736293956
$ws.Range("I10").Value = 0.4387957959343041
$ws.Range("J10").Value = 0.03298154123656971
$ws.Range("M10").Value = 0.70251001423415
$ws.Range("B11").Value = 1.281650434587448
$ws.Range("C11").Value = 0.3582313320918615
$ws.Range("E11").Value = 0.4829092345509238
$ws.Range("F11").Value = 2.142075605049072
$ws.Range("G11").Value = 0.5586165367839016
$ws.Range("H11").Value = 0.6320227970032022
$ws.Range("I11").Value = 0.4348228994163676
$ws.Range("J11").Value = 0.03244221498611655
$ws.Range("M11").Value = 0.7327391807038737
$ws.Range("B12").Value = 1.312488792671502
$ws.Range("C12").Value = 0.3671184847376594
$ws.Range("E12").Value = 0.484240595068826
$ws.Range("F12").Value = 2.148665991630224
$ws.Range("G12").Value = 0.5608963330929697
$ws.Range("H12").Value = 0.6315872744150681
$ws.Range("I12").Value = 0.4333912683193653
$ws.Range("J12").Value = 0.03224240001847178
$ws.Range("M12").Value = 0.7442097062976529
$ws.Range("B13").Value = 1.305847728232038
$ws.Range("C13").Value = 0.365204794226031
$ws.Range("E13").Value = 0.4839528612040098
$ws.Range("F13").Value = 2.147237019594144
$ws.Range("G13").Value = 0.5604007434939007
$ws.Range("H13").Value = 0.6316775883430239
$ws.Range("I13").Value = 0.4336963487479366
$ws.Range("J13").Value = 0.03228523701104269
$ws.Range("M13").Value = 0.741738287436263
$ws.Range("B14").Value = 1.284187761527676
$ws.Range("C14").Value = 0.3589626261565968
$ws.Range("E14").Value = 0.4830183149496818
$ws.Range("F14").Value = 2.142613468956185
$ws.Range("G14").Value = 0.5588020281297617
$ws.Range("H14").Value = 0.6319853970750131
$ws.Range("I14").Value = 0.4347036558804156
$ws.Range("J14").Value = 0.03242568746480767
$ws.Range("M14").Value = 0.7336824005136862
$ws.Range("B15").Value = 1.270918874693564
$ws.Range("C15").Value = 0.3551381914305125
$ws.Range("E15").Value = 0.4824488123637067
$ws.Range("F15").Value = 2.139809547949469
$ws.Range("G15").Value = 0.5578362040399583
$ws.Range("H15").Value = 0.6321841325749062
$ws.Range("I15").Value = 0.4353301592113183
$ws.Range("J15").Value = 0.03251229314545778
$ws.Range("M15").Value = 0.7287509754300885
$ws.Range("B16").Value = 1.194857367668874
$ws.Range("C16").Value = 0.3332077524809449
$ws.Range("E16").Value = 0.4792302501948669
$ws.Range("F16").Value = 2.124172024573411
$ws.Range("G16").Value = 0.5525065653730081
$ws.Range("H16").Value = 0.6334789478233631
$ws.Range("I16").Value = 0.4390655955455216
$ws.Range("J16").Value = 0.03301740434741696
$ws.Range("M16").Value = 0.7005377569743558
$ws.Range("B17").Value = 1.148176335376093
$ws.Range("C17").Value = 0.3197415703390334
$ws.Range("E17").Value = 0.4772960072214687
$ws.Range("F17").Value = 2.114962446554301
$ws.Range("G17").Value = 0.5494191391426568
$ws.Range("H17").Value = 0.6344120992607429
$ws.Range("I17").Value = 0.4414863056823037
$ws.Range("J17").Value = 0.03333511774797993
$ws.Range("M17").Value = 0.6832718461108556
$ws.Range("B18").Value = 1.121319640476202
$ws.Range("C18").Value = 0.3119915833690357
$ws.Range("E18").Value = 0.4761982949919314
$ws.Range("F18").Value = 2.109806008460808
$ws.Range("G18").Value = 0.5477099151328133
$ws.Range("H18").Value = 0.6349997115511172
$ws.Range("I18").Value = 0.4429259254162226
$ws.Range("J18").Value = 0.03352073425957247
$ws.Range("M18").Value = 0.6733564752422865
$ws.Range("B19").Value = 1.112225247944366
$ws.Range("C19").Value = 0.3093667840129513
$ws.Range("E19").Value = 0.4758291750406727
$ws.Range("F19").Value = 2.108084260364336
$ws.Range("G19").Value = 0.5471426033925297
$ws.Range("H19").Value = 0.635207394854703
$ws.Range("I19").Value = 0.4434214623664587
$ws.Range("J19").Value = 0.03358407419282816
$ws.Range("M19").Value = 0.6700019764664944
$ws.Range("B20").Value = 1.153146341368426
$ws.Range("C20").Value = 0.3211755445447864
$ws.Range("E20").Value = 0.4775003781436737
$ws.Range("F20").Value = 2.115928256961737
$ws.Range("G20").Value = 0.5497409020229043
$ws.Range("H20").Value = 0.6343074940378131
$ws.Range("I20").Value = 0.441223718704002
$ws.Range("J20").Value = 0.03330099879560589
$ws.Range("M20").Value = 0.6851082272414288
$ws.Range("B21").Value = 1.29055013944452
$ws.Range("C21").Value = 0.3607962954438335
$ws.Range("E21").Value = 0.4832922025356652
$ws.Range("F21").Value = 2.143965652507788
$ws.Range("G21").Value = 0.559268807738647
$ws.Range("H21").Value = 0.6318928610975547
$ws.Range("I21").Value = 0.4344058049204307
$ws.Range("J21").Value = 0.03238431374533235
$ws.Range("M21").Value = 0.7360479771245423
$ws.Range("B22").Value = 1.380283606292551
$ws.Range("C22").Value = 0.3866493389600691
$ws.Range("E22").Value = 0.4872088791728331
$ws.Range("F22").Value = 2.163548312793367
$ws.Range("G22").Value = 0.5660962273061614
$ws.Range("H22").Value = 0.6307706599800582
$ws.Range("I22").Value = 0.4303745892997597
$ws.Range("J22").Value = 0.03181095752110874
$ws.Range("M22").Value = 0.7694762883729567
$ws.Range("B23").Value = 1.332397638717282
$ws.Range("C23").Value = 0.3728548944382055
$ws.Range("E23").Value = 0.4851064785922929
$ws.Range("F23").Value = 2.152981226209832
$ws.Range("G23").Value = 0.5623970077281939
$ws.Range("H23").Value = 0.6313277516210292
$ws.Range("I23").Value = 0.4324870936394483
$ws.Range("J23").Value = 0.03211460560410462
$ws.Range("M23").Value = 0.7516225990834187
$ws.Range("B24").Value = 1.150899461249196
$ws.Range("C24").Value = 0.320527270036564
$ws.Range("E24").Value = 0.4774079374705806
$ws.Range("F24").Value = 2.115491183425988
$ws.Range("G24").Value = 0.5495952282085597
$ws.Range("H24").Value = 0.6343546268106195
$ws.Range("I24").Value = 0.4413422850381785
$ws.Range("J24").Value = 0.03331641475865288
$ws.Range("M24").Value = 0.6842779649520736
$ws.Range("B25").Value = 0.9548594282517229
$ws.Range("C25").Value = 0.2639079056741593
$ws.Range("E25").Value = 0.4696704338484921
$ws.Range("F25").Value = 2.080426350308073
$ws.Range("G25").Value = 0.5383274623210781
$ws.Range("H25").Value = 0.6395635181086305
$ws.Range("I25").Value = 0.452699976984313
$ws.Range("J25").Value = 0.03472264513235546
$ws.Range("M25").Value = 0.612230378514667
